$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Refresh crawl timestamp (column O) for every data row (2-33).
# -----------------------------------------------------------------------
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-14 20:57:05"
}

# -----------------------------------------------------------------------
# 2. Re-order the product rows (the crawl re-ran and produced the items
#    in a different order). Rows are shuffled using a scratch row
#    (row 100) as a temporary holding area so cell types/content are
#    preserved exactly as copied (rather than retyped, which Excel would
#    otherwise try to reinterpret, e.g. turning "3.95" into a number).
# -----------------------------------------------------------------------

# Rows 11-12 swap: new11 = old12, new12 = old11
$ws.Range("A11:O11").Copy($ws.Range("A100:O100"))
$ws.Range("A12:O12").Copy($ws.Range("A11:O11"))
$ws.Range("A100:O100").Copy($ws.Range("A12:O12"))
$ws.Range("A100:O100").Clear()
# Restore blank cell lost by the copy-over-nonblank limitation.
$ws.Range("E11").ClearContents()

# Rows 19-21 rotate: new19 = old21, new20 = old19, new21 = old20
$ws.Range("A19:O19").Copy($ws.Range("A100:O100"))
$ws.Range("A21:O21").Copy($ws.Range("A19:O19"))
$ws.Range("A20:O20").Copy($ws.Range("A21:O21"))
$ws.Range("A100:O100").Copy($ws.Range("A20:O20"))
$ws.Range("A100:O100").Clear()
$ws.Range("E19").ClearContents()

# Rows 25-27 rotate: new25 = old27, new26 = old25, new27 = old26
$ws.Range("A25:O25").Copy($ws.Range("A100:O100"))
$ws.Range("A27:O27").Copy($ws.Range("A25:O25"))
$ws.Range("A26:O26").Copy($ws.Range("A27:O27"))
$ws.Range("A100:O100").Copy($ws.Range("A26:O26"))
$ws.Range("A100:O100").Clear()
$ws.Range("E26").ClearContents()
